$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.035434270339478
$ws.Cells.Item(2, 4).Value = 1.03369087082837
$ws.Cells.Item(2, 5).Value = 1.043521688405821
$ws.Cells.Item(2, 6).Value = 1.052325948720622
$ws.Cells.Item(2, 9).Value = 1.032945189209992
$ws.Cells.Item(2, 10).Value = 1.040548066182423
$ws.Cells.Item(2, 11).Value = 1.036492492875783
$ws.Cells.Item(2, 12).Value = 1.04629531870936
$ws.Cells.Item(2, 13).Value = 1.055074991578402
$ws.Cells.Item(2, 14).Value = 1.017192570817582

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.037196827573088
$ws.Cells.Item(3, 4).Value = 1.034129597227759
$ws.Cells.Item(3, 5).Value = 1.045080494862234
$ws.Cells.Item(3, 6).Value = 1.053952840104193
$ws.Cells.Item(3, 9).Value = 1.033032541169069
$ws.Cells.Item(3, 10).Value = 1.041950070006049
$ws.Cells.Item(3, 11).Value = 1.036741240582125
$ws.Cells.Item(3, 12).Value = 1.047663223591214
$ws.Cells.Item(3, 13).Value = 1.056512617831267
$ws.Cells.Item(3, 14).Value = 1.017681221255893

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.038334867433945
$ws.Cells.Item(4, 4).Value = 1.034412835668525
$ws.Cells.Item(4, 5).Value = 1.046086972311852
$ws.Cells.Item(4, 6).Value = 1.055003019729354
$ws.Cells.Item(4, 9).Value = 1.033087175249729
$ws.Cells.Item(4, 10).Value = 1.042854568334637
$ws.Cells.Item(4, 11).Value = 1.036900820779807
$ws.Cells.Item(4, 12).Value = 1.048545692905292
$ws.Cells.Item(4, 13).Value = 1.057439858918786
$ws.Cells.Item(4, 14).Value = 1.017995881716732

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.038812726966742
$ws.Cells.Item(5, 4).Value = 1.034531753038344
$ws.Cells.Item(5, 5).Value = 1.046509586566416
$ws.Cells.Item(5, 6).Value = 1.055443921971726
$ws.Cells.Item(5, 9).Value = 1.033109691271158
$ws.Cells.Item(5, 10).Value = 1.043234186462089
$ws.Cells.Item(5, 11).Value = 1.036967578214165
$ws.Cells.Item(5, 12).Value = 1.048916057831095
$ws.Cells.Item(5, 13).Value = 1.057828963970866
$ws.Cells.Item(5, 14).Value = 1.018127802744199

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.038892928542982
$ws.Cells.Item(6, 4).Value = 1.034551710600387
$ws.Cells.Item(6, 5).Value = 1.046580515858118
$ws.Cells.Item(6, 6).Value = 1.055517916952532
$ws.Cells.Item(6, 9).Value = 1.033113445294711
$ws.Cells.Item(6, 10).Value = 1.043297889229946
$ws.Cells.Item(6, 11).Value = 1.036978767706788
$ws.Cells.Item(6, 12).Value = 1.04897820738402
$ws.Cells.Item(6, 13).Value = 1.057894255271886
$ws.Cells.Item(6, 14).Value = 1.018149931734307

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.038341254845494
$ws.Cells.Item(7, 4).Value = 1.034414425262879
$ws.Cells.Item(7, 5).Value = 1.04609262128756
$ws.Cells.Item(7, 6).Value = 1.055008913401431
$ws.Cells.Item(7, 9).Value = 1.033087477886468
$ws.Cells.Item(7, 10).Value = 1.042859643282012
$ws.Cells.Item(7, 11).Value = 1.036901714093554
$ws.Cells.Item(7, 12).Value = 1.04855064418117
$ws.Cells.Item(7, 13).Value = 1.057445060922897
$ws.Cells.Item(7, 14).Value = 1.017997645869982

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.036030451513626
$ws.Cells.Item(8, 4).Value = 1.033839273690455
$ws.Cells.Item(8, 5).Value = 1.044048951693083
$ws.Cells.Item(8, 6).Value = 1.052876294895634
$ws.Cells.Item(8, 9).Value = 1.032975101215413
$ws.Cells.Item(8, 10).Value = 1.041022444249276
$ws.Cells.Item(8, 11).Value = 1.036576842915293
$ws.Cells.Item(8, 12).Value = 1.046758165181541
$ws.Cells.Item(8, 13).Value = 1.055561471714006
$ws.Cells.Item(8, 14).Value = 1.017358031435108

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.031939074479626
$ws.Cells.Item(9, 4).Value = 1.032820880923414
$ws.Cells.Item(9, 5).Value = 1.040430561481363
$ws.Cells.Item(9, 6).Value = 1.049098450987253
$ws.Cells.Item(9, 9).Value = 1.032762616498003
$ws.Cells.Item(9, 10).Value = 1.03776393087929
$ws.Cells.Item(9, 11).Value = 1.035993864392361
$ws.Cells.Item(9, 12).Value = 1.043578753614481
$ws.Cells.Item(9, 13).Value = 1.052218883744749
$ws.Cells.Item(9, 14).Value = 1.016219065580962

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.02919748458759
$ws.Cells.Item(10, 4).Value = 1.032138724530178
$ws.Cells.Item(10, 5).Value = 1.038006005078891
$ws.Cells.Item(10, 6).Value = 1.046565754991037
$ws.Cells.Item(10, 9).Value = 1.032611235019801
$ws.Cells.Item(10, 10).Value = 1.035576652469598
$ws.Cells.Item(10, 11).Value = 1.035598175952241
$ws.Cells.Item(10, 12).Value = 1.041444451710952
$ws.Cells.Item(10, 13).Value = 1.049974023300615
$ws.Cells.Item(10, 14).Value = 1.015451521655385

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.028006815935777
$ws.Cells.Item(11, 4).Value = 1.031842591373054
$ws.Cells.Item(11, 5).Value = 1.036953062389154
$ws.Cells.Item(11, 6).Value = 1.045465548047732
$ws.Cells.Item(11, 9).Value = 1.03254337757278
$ws.Cells.Item(11, 10).Value = 1.034625831321246
$ws.Cells.Item(11, 11).Value = 1.035425175869674
$ws.Cells.Item(11, 12).Value = 1.040516636375775
$ws.Cells.Item(11, 13).Value = 1.048997912322495
$ws.Cells.Item(11, 14).Value = 1.015117159059255

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.027563998401609
$ws.Cells.Item(12, 4).Value = 1.031732481955445
$ws.Cells.Item(12, 5).Value = 1.036561473345419
$ws.Cells.Item(12, 6).Value = 1.045056336719752
$ws.Cells.Item(12, 9).Value = 1.032517825337996
$ws.Cells.Item(12, 10).Value = 1.03427208172782
$ws.Cells.Item(12, 11).Value = 1.03536066630855
$ws.Cells.Item(12, 12).Value = 1.04017144294781
$ws.Cells.Item(12, 13).Value = 1.048634716010404
$ws.Cells.Item(12, 14).Value = 1.014992654995565

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.027659009454146
$ws.Cells.Item(13, 4).Value = 1.031756105868167
$ws.Cells.Item(13, 5).Value = 1.036645492445112
$ws.Cells.Item(13, 6).Value = 1.045144138860896
$ws.Cells.Item(13, 9).Value = 1.032523322082516
$ws.Cells.Item(13, 10).Value = 1.034347988341905
$ws.Cells.Item(13, 11).Value = 1.035374515117966
$ws.Cells.Item(13, 12).Value = 1.040245513734996
$ws.Cells.Item(13, 13).Value = 1.048712651364154
$ws.Cells.Item(13, 14).Value = 1.015019375505401

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.027970223828762
$ws.Cells.Item(14, 4).Value = 1.031833491980787
$ws.Cells.Item(14, 5).Value = 1.036920703377269
$ws.Cells.Item(14, 6).Value = 1.045431733725094
$ws.Cells.Item(14, 9).Value = 1.03254127249594
$ws.Cells.Item(14, 10).Value = 1.034596602026187
$ws.Cells.Item(14, 11).Value = 1.035419848579175
$ws.Cells.Item(14, 12).Value = 1.040488114111587
$ws.Cells.Item(14, 13).Value = 1.048967903235413
$ws.Cells.Item(14, 14).Value = 1.015106873803183

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.028161899829406
$ws.Cells.Item(15, 4).Value = 1.031881157258411
$ws.Cells.Item(15, 5).Value = 1.037090205975958
$ws.Cells.Item(15, 6).Value = 1.045608857688032
$ws.Cells.Item(15, 9).Value = 1.032552286357553
$ws.Cells.Item(15, 10).Value = 1.034749704863978
$ws.Cells.Item(15, 11).Value = 1.035447746955084
$ws.Cells.Item(15, 12).Value = 1.040637513403958
$ws.Cells.Item(15, 13).Value = 1.049125089080877
$ws.Cells.Item(15, 14).Value = 1.015160743580069

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.029276429182645
$ws.Cells.Item(16, 4).Value = 1.032158362096536
$ws.Cells.Item(16, 5).Value = 1.038075819006261
$ws.Cells.Item(16, 6).Value = 1.046638696417366
$ws.Cells.Item(16, 9).Value = 1.032615689848718
$ws.Cells.Item(16, 10).Value = 1.03563967586177
$ws.Cells.Item(16, 11).Value = 1.035609622368353
$ws.Cells.Item(16, 12).Value = 1.041505949752652
$ws.Cells.Item(16, 13).Value = 1.050038717637092
$ws.Cells.Item(16, 14).Value = 1.015473669423649

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.029974583013321
$ws.Cells.Item(17, 4).Value = 1.032332043994542
$ws.Cells.Item(17, 5).Value = 1.038693229586195
$ws.Cells.Item(17, 6).Value = 1.047283730864842
$ws.Cells.Item(17, 9).Value = 1.032654843082645
$ws.Cells.Item(17, 10).Value = 1.036196926581372
$ws.Cells.Item(17, 11).Value = 1.03571071724533
$ws.Cells.Item(17, 12).Value = 1.042049710738166
$ws.Cells.Item(17, 13).Value = 1.050610713114638
$ws.Cells.Item(17, 14).Value = 1.015669417779463

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.030381463793026
$ws.Cells.Item(18, 4).Value = 1.032433276798916
$ws.Cells.Item(18, 5).Value = 1.039053057212056
$ws.Cells.Item(18, 6).Value = 1.047659628677356
$ws.Cells.Item(18, 9).Value = 1.032677457799857
$ws.Cells.Item(18, 10).Value = 1.036521604000939
$ws.Cells.Item(18, 11).Value = 1.035769523561727
$ws.Cells.Item(18, 12).Value = 1.042366526201189
$ws.Cells.Item(18, 13).Value = 1.050943956471512
$ws.Cells.Item(18, 14).Value = 1.015783400775957

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.030520142405105
$ws.Cells.Item(19, 4).Value = 1.032467782170206
$ws.Cells.Item(19, 5).Value = 1.039175699103591
$ws.Cells.Item(19, 6).Value = 1.047787742912237
$ws.Cells.Item(19, 9).Value = 1.032685131064524
$ws.Cells.Item(19, 10).Value = 1.036632250397652
$ws.Cells.Item(19, 11).Value = 1.035789547739303
$ws.Cells.Item(19, 12).Value = 1.042474492932317
$ws.Cells.Item(19, 13).Value = 1.051057517773256
$ws.Cells.Item(19, 14).Value = 1.015822233319549

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.029899713068031
$ws.Cells.Item(20, 4).Value = 1.032313417104346
$ws.Cells.Item(20, 5).Value = 1.03862701815222
$ws.Cells.Item(20, 6).Value = 1.04721456003349
$ws.Cells.Item(20, 9).Value = 1.032650665347487
$ws.Cells.Item(20, 10).Value = 1.036137175932497
$ws.Cells.Item(20, 11).Value = 1.035699887331889
$ws.Cells.Item(20, 12).Value = 1.041991406744922
$ws.Cells.Item(20, 13).Value = 1.050549384049782
$ws.Cells.Item(20, 14).Value = 1.015648435904005

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.027878594285656
$ws.Cells.Item(21, 4).Value = 1.031810706793355
$ws.Cells.Item(21, 5).Value = 1.036839673963553
$ws.Cells.Item(21, 6).Value = 1.045347059382663
$ws.Cells.Item(21, 9).Value = 1.03253599612574
$ws.Cells.Item(21, 10).Value = 1.034523407404736
$ws.Cells.Item(21, 11).Value = 1.035406505892391
$ws.Cells.Item(21, 12).Value = 1.04041668995551
$ws.Cells.Item(21, 13).Value = 1.048892755290202
$ws.Cells.Item(21, 14).Value = 1.015081116246134

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.026604644434086
$ws.Cells.Item(22, 4).Value = 1.031493983181272
$ws.Cells.Item(22, 5).Value = 1.035713119686494
$ws.Cells.Item(22, 6).Value = 1.044169724198714
$ws.Cells.Item(22, 9).Value = 1.032461891549654
$ws.Cells.Item(22, 10).Value = 1.033505448912615
$ws.Cells.Item(22, 11).Value = 1.035220601262663
$ws.Cells.Item(22, 12).Value = 1.03942334761012
$ws.Cells.Item(22, 13).Value = 1.047847542815665
$ws.Cells.Item(22, 14).Value = 1.014722641776887

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.027280298009889
$ws.Cells.Item(23, 4).Value = 1.031661945602234
$ws.Cells.Item(23, 5).Value = 1.036310595804439
$ws.Cells.Item(23, 6).Value = 1.044794156664648
$ws.Cells.Item(23, 9).Value = 1.032501366133955
$ws.Cells.Item(23, 10).Value = 1.034045407249994
$ws.Cells.Item(23, 11).Value = 1.035319289569375
$ws.Cells.Item(23, 12).Value = 1.039950250253424
$ws.Cells.Item(23, 13).Value = 1.048401977814483
$ws.Cells.Item(23, 14).Value = 1.01491284602706

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.029933544625557
$ws.Cells.Item(24, 4).Value = 1.032321834019883
$ws.Cells.Item(24, 5).Value = 1.038656937168404
$ws.Cells.Item(24, 6).Value = 1.04724581640485
$ws.Cells.Item(24, 9).Value = 1.032652553774431
$ws.Cells.Item(24, 10).Value = 1.036164175783279
$ws.Cells.Item(24, 11).Value = 1.035704781400131
$ws.Cells.Item(24, 12).Value = 1.042017752894628
$ws.Cells.Item(24, 13).Value = 1.050577097225272
$ws.Cells.Item(24, 14).Value = 1.015657917308899

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.032999193218948
$ws.Cells.Item(25, 4).Value = 1.033084733688609
$ws.Cells.Item(25, 5).Value = 1.041368114414363
$ws.Cells.Item(25, 6).Value = 1.050077549419802
$ws.Cells.Item(25, 9).Value = 1.032819262844663
$ws.Cells.Item(25, 10).Value = 1.038608913530145
$ws.Cells.Item(25, 11).Value = 1.036145820804346
$ws.Cells.Item(25, 12).Value = 1.044403246905675
$ws.Cells.Item(25, 13).Value = 1.053085874245875
$ws.Cells.Item(25, 14).Value = 1.016514948707568
